# Timeline.xlsx update: "update timeline, upload interview responses"
#
# Adds four new task notes to the existing timeline grid, grows the row
# heights that now need to wrap the extra text, and updates the active
# selection to where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New task notes -------------------------------------------------
# Backend row (row 3): note under the 2023-03-18 column (H)
$ws.Range("H3").Value = "Consolidate all interview responses for backend"

# Modelling row (row 4): notes under 2023-03-22 (L) and 2023-03-24 (N)
$ws.Range("L4").Value = "Get model and data from backend"
$ws.Range("N4").Value = "Implement skeleton interactive viz based on model and actual data"

# Consultation row (row 5): note under 2023-03-24 (N)
$ws.Range("N5").Value = "Integrate frontend, backend during meeting"

# Consultation row in the second block (row 10): same note, under 2023-03-24 (N)
$ws.Range("N10").Value = "Integrate frontend, backend during meeting"

# --- Row heights grown to fit the new wrapped text -------------------
$ws.Rows.Item(3).RowHeight = 85
$ws.Rows.Item(4).RowHeight = 119
$ws.Rows.Item(5).RowHeight = 86
$ws.Rows.Item(10).RowHeight = 86

# --- Update the active selection -------------------------------------
$ws.Range("S5").Select() | Out-Null
